# Applies the FlashScore odds-update edit described in the commit diff.
# Only numeric odds/statistics cells change; no structural/style edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("L9").Value = 2.9
$ws.Range("U9").Value = 1.7
$ws.Range("V9").Value = 1.93
$ws.Range("AD9").Value = 5.9
$ws.Range("AE9").Value = 13
$ws.Range("AG9").Value = 400
$ws.Range("AH9").Value = 8
$ws.Range("AN9").Value = 4.9
$ws.Range("AO9").Value = 17.5
$ws.Range("AP9").Value = 24
$ws.Range("AT9").Value = 2.42
$ws.Range("AU9").Value = 6.7
$ws.Range("AV9").Value = 60
$ws.Range("AW9").Value = 4.2
$ws.Range("AX9").Value = 12
$ws.Range("AY9").Value = 19.5

# Row 19
$ws.Range("G19").Value = 1.53
$ws.Range("I19").Value = 5.25
$ws.Range("J19").Value = 2
$ws.Range("L19").Value = 5
$ws.Range("M19").Value = 1.01
$ws.Range("N19").Value = 23
$ws.Range("Q19").Value = 1.4
$ws.Range("R19").Value = 2.88
$ws.Range("U19").Value = 1.53
$ws.Range("V19").Value = 2.38
$ws.Range("W19").Value = 11
$ws.Range("AB19").Value = 19
$ws.Range("AD19").Value = 10
$ws.Range("AE19").Value = 15
$ws.Range("AF19").Value = 41
$ws.Range("AI19").Value = 34
$ws.Range("AS19").Value = 81
$ws.Range("AU19").Value = 7.5
$ws.Range("AW19").Value = 7.5
$ws.Range("AZ19").Value = 81
$ws.Range("BA19").Value = 81
$ws.Range("BB19").Value = 126
$ws.Range("BC19").Value = 301

# Row 20
$ws.Range("G20").Value = 1.6
$ws.Range("H20").Value = 4.33
$ws.Range("I20").Value = 5.25
$ws.Range("J20").Value = 2.1
$ws.Range("K20").Value = 2.4
$ws.Range("O20").Value = 1.18
$ws.Range("P20").Value = 4.5
$ws.Range("Q20").Value = 1.62
$ws.Range("R20").Value = 2.25
$ws.Range("Z20").Value = 12
$ws.Range("AC20").Value = 15
$ws.Range("AD20").Value = 8
$ws.Range("AJ20").Value = 17
$ws.Range("BA20").Value = 101

# Row 21
$ws.Range("M21").Value = 1.02
$ws.Range("N21").Value = 19
$ws.Range("Q21").Value = 1.5
$ws.Range("R21").Value = 2.5

# Row 22
$ws.Range("G22").Value = 1.8
$ws.Range("I22").Value = 3.8
$ws.Range("J22").Value = 2.3
$ws.Range("N22").Value = 21
$ws.Range("AD22").Value = 9
$ws.Range("AG22").Value = 81
$ws.Range("AI22").Value = 23

# Row 23
$ws.Range("G23").Value = 2.45
$ws.Range("I23").Value = 2.8
$ws.Range("J23").Value = 3
$ws.Range("N23").Value = 13
$ws.Range("Q23").Value = 1.73
$ws.Range("R23").Value = 2.08
$ws.Range("Y23").Value = 9.5
$ws.Range("Z23").Value = 23
$ws.Range("AB23").Value = 23
$ws.Range("AC23").Value = 13
$ws.Range("AH23").Value = 11
$ws.Range("AK23").Value = 29
$ws.Range("AN23").Value = 4.5
$ws.Range("AW23").Value = 5
$ws.Range("AZ23").Value = 51
$ws.Range("BC23").Value = 451

# Row 24
$ws.Range("G24").Value = 2.1
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = 3
$ws.Range("J24").Value = 2.6
$ws.Range("L24").Value = 3.25
$ws.Range("O24").Value = 1.1
$ws.Range("P24").Value = 7
$ws.Range("W24").Value = 15
$ws.Range("X24").Value = 15
$ws.Range("Y24").Value = 9.5
$ws.Range("Z24").Value = 21
$ws.Range("AA24").Value = 15
$ws.Range("AF24").Value = 26
$ws.Range("AH24").Value = 17
$ws.Range("AI24").Value = 21
$ws.Range("AK24").Value = 34
$ws.Range("AL24").Value = 21
$ws.Range("AN24").Value = 4.75
$ws.Range("AO24").Value = 11
$ws.Range("AQ24").Value = 34
$ws.Range("AW24").Value = 5.5
$ws.Range("AX24").Value = 15
$ws.Range("AY24").Value = 17
$ws.Range("AZ24").Value = 41
$ws.Range("BD24").Value = 176

# Row 25
$ws.Range("G25").Value = 1.53
$ws.Range("H25").Value = 4.33
$ws.Range("I25").Value = 6
$ws.Range("J25").Value = 2
$ws.Range("X25").Value = 9
$ws.Range("AJ25").Value = 19
$ws.Range("AN25").Value = 3.75
$ws.Range("AR25").Value = 34
$ws.Range("AW25").Value = 7.5
$ws.Range("BB25").Value = 151

# Row 26
$ws.Range("O26").Value = 1.25
$ws.Range("P26").Value = 3.75

# Row 27
$ws.Range("G27").Value = 2.25
$ws.Range("I27").Value = 2.8
$ws.Range("L27").Value = 3
$ws.Range("N27").Value = 29
$ws.Range("W27").Value = 19
$ws.Range("X27").Value = 19
$ws.Range("Z27").Value = 26
$ws.Range("AF27").Value = 21
$ws.Range("AN27").Value = 5.5
$ws.Range("AP27").Value = 13

# Row 28
$ws.Range("G28").Value = 1.91
$ws.Range("H28").Value = 3.9
$ws.Range("I28").Value = 3.5
$ws.Range("J28").Value = 2.4
$ws.Range("L28").Value = 3.6
$ws.Range("N28").Value = 21
$ws.Range("AH28").Value = 17
$ws.Range("AJ28").Value = 12
$ws.Range("AO28").Value = 10

# Row 29
$ws.Range("G29").Value = 1.38
$ws.Range("H29").Value = 5.75
$ws.Range("I29").Value = 6.5
$ws.Range("J29").Value = 1.8
$ws.Range("K29").Value = 2.88
$ws.Range("L29").Value = 6
$ws.Range("N29").Value = 29
$ws.Range("U29").Value = 1.53
$ws.Range("V29").Value = 2.38
$ws.Range("W29").Value = 13
$ws.Range("Y29").Value = 9
$ws.Range("Z29").Value = 11
$ws.Range("AC29").Value = 29
$ws.Range("AD29").Value = 12
$ws.Range("AE29").Value = 17
$ws.Range("AG29").Value = 126
$ws.Range("AO29").Value = 6.5
$ws.Range("AQ29").Value = 15
$ws.Range("AW29").Value = 9
$ws.Range("AZ29").Value = 101

# Row 30
$ws.Range("O30").Value = 1.1
$ws.Range("P30").Value = 7

# Row 31
$ws.Range("N31").Value = 17

# Row 32
$ws.Range("M32").Value = 1.05
$ws.Range("N32").Value = 11
$ws.Range("Q32").Value = 1.85
$ws.Range("R32").Value = 2
